# Update attendee/view counts (column F) across the 展览 / 本地生活 / 全部类型
# sheets to match the newly scraped numbers (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

# ---- 展览 ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 12682
$ws.Range("F3").Value  = 7095
$ws.Range("F6").Value  = 445
$ws.Range("F10").Value = 995
$ws.Range("F11").Value = 140
$ws.Range("F12").Value = 348
$ws.Range("F13").Value = 997
$ws.Range("F15").Value = 1012
$ws.Range("F18").Value = 365
$ws.Range("F24").Value = 364
$ws.Range("F25").Value = 5208
$ws.Range("F27").Value = 1415
$ws.Range("F29").Value = 1290
$ws.Range("F30").Value = 1290
$ws.Range("F31").Value = 46
$ws.Range("F32").Value = 20
$ws.Range("F33").Value = 1331
$ws.Range("F34").Value = 3
$ws.Range("F38").Value = 3719

# ---- 本地生活 ------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9255
$ws.Range("F4").Value = 1981

# ---- 全部类型 (aggregated view of all the other sheets) -------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9255
$ws.Range("F4").Value  = 1981
$ws.Range("F6").Value  = 12682
$ws.Range("F7").Value  = 7095
$ws.Range("F10").Value = 445
$ws.Range("F13").Value = 140
$ws.Range("F14").Value = 348
$ws.Range("F15").Value = 997
$ws.Range("F17").Value = 1012
$ws.Range("F20").Value = 365
$ws.Range("F29").Value = 364
$ws.Range("F30").Value = 5208
$ws.Range("F32").Value = 1415
$ws.Range("F37").Value = 1290
$ws.Range("F38").Value = 1290
$ws.Range("F39").Value = 1331
$ws.Range("F40").Value = 3
$ws.Range("F48").Value = 3719
